{"js": "// Apply the benchmark-table value updates described by the commit.\n// The document body contains a single one-column table; each row holds\n// one paragraph/run of text. A handful of rows get their numeric value\n// replaced, and three rows that previously packed 10 tab-separated\n// values into one run get collapsed down to a single value.\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// rowIndex -> new cell text (0-based row indices)\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"80\",\n  5: \"0.00075\",\n  6: \"0.00026\",\n  7: \"0.00010\",\n  8: \"0.00034\",\n  9: \"0.00045\",\n  10: \"0.00052\",\n  11: \"0.02106\",\n  43: \"99.97\",\n  44: \"0.02\",\n  45: \"70\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const cell = table.getCell(Number(rowIndex), 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-table value updates described by the commit.\n# The document body contains a single one-column table; each row holds\n# one paragraph/run of text. A handful of rows get their numeric value\n# replaced, and three rows that previously packed 10 tab-separated\n# values into one run get collapsed down to a single value.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# 1-based row index -> new cell text\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"80\"\n    6  = \"0.00075\"\n    7  = \"0.00026\"\n    8  = \"0.00010\"\n    9  = \"0.00034\"\n    10 = \"0.00045\"\n    11 = \"0.00052\"\n    12 = \"0.02106\"\n    44 = \"99.97\"\n    45 = \"0.02\"\n    46 = \"70\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $t.Cell($rowIndex, 1).Range.Text = $updates[$rowIndex]\n}\n"}
